$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1) Split the existing "Key Class Design..." entry (F18) so that "추천로직"
#    becomes its own run with an explicit 맑은 고딕 / 10pt font, matching the
#    rich-text run layout introduced upstream.
# ---------------------------------------------------------------------------
$f18Chars = $ws.Range("F18").Characters(51, 4)
$f18Chars.Font.Name = "맑은 고딕"
$f18Chars.Font.Size = 10

# ---------------------------------------------------------------------------
# 2) Fill in row 19 (2019-11-12 PSP log entry).
# ---------------------------------------------------------------------------
$ws.Range("A19").Value = 43781
$ws.Range("A19").NumberFormat = 'm"월"\ d"일";@'
$ws.Range("B19").Value = 0.79166666666666663
$ws.Range("C19").Value = 0.91666666666666663
$ws.Range("D19").Value = 30
$ws.Range("E19").Value = 150
$ws.Range("F19").Value = "추천로직 논의"

$f19Chars = $ws.Range("F19").Characters(5, 3)
$f19Chars.Font.Name = "돋움"
$f19Chars.Font.Size = 10

$ws.Rows.Item(19).RowHeight = 13

# ---------------------------------------------------------------------------
# 3) Fill in row 20 (2019-11-13 PSP log entry).
# ---------------------------------------------------------------------------
$ws.Range("A20").Value = 43782
$ws.Range("A20").NumberFormat = 'm"월"\ d"일";@'
$ws.Range("B20").Value = 0
$ws.Range("C20").Value = 0.083333333333333329
$ws.Range("D20").Value = 0
$ws.Range("E20").Value = 120
$ws.Range("F20").Value = "Key Class Design, Table Design 수정 및 용어 통일"

$f20Chars = $ws.Range("F20").Characters(34, 8)
$f20Chars.Font.Name = "돋움"
$f20Chars.Font.Size = 10

# ---------------------------------------------------------------------------
# 4) Move the active selection to D20, matching where editing left off.
# ---------------------------------------------------------------------------
$ws.Range("D20").Select()
